$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enter the two values in the order they were authored so the shared-string
# table indices come out as "hkshsfh" = 0, "isf" = 1 (matches the target
# sharedStrings.xml / sheet1.xml <v> indices).
$ws.Range("E8").Value = "hkshsfh"
$ws.Range("F6").Value = "isf"

# Leave the selection on F6, as captured in the saved sheetView.
$ws.Range("F6").Select()
